$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LeetCode Algorithm")
$ws.Range("A2481:H2482").UnMerge()
$ws.Range("2481:2631").Insert()

$ws.Range("A2632:H2632").Merge()
$ws.Range("A2633:H2633").Merge()

$ws.Range("A2481").Value = 2792
$ws.Range("B2").Copy()
$ws.Range("B2481").PasteSpecial(-4122)
$ws.Hyperlinks.Add($ws.Range("B2481"), "https://leetcode.com/problems/count-nodes-that-are-great-enough", "", "", "https://leetcode.com/problems/count-nodes-that-are-great-enough")
$ws.Range("B2481").Value = "Count Nodes That Are Great Enough"
$ws.Range("C2481").Value = 4
$ws.Range("D2481").Value = "C"
$ws.Range("F2481").Value = "Hard"
$ws.Range("G2481").Value = "Tree"
$ws.Range("H2481").Value = "Track the k smallest nodes in the subtree"
$ws.Range("I2481").Value = "N * LOG(K)"
